$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.146286725997925
$ws.Range("B1").Value = 3.700870752334595
$ws.Range("C1").Value = 4.405519962310791
$ws.Range("D1").Value = 1.851027607917786
$ws.Range("E1").Value = 1.285533785820007
